$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-04-17 Wednesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-04-18 Thursday", 2) | Out-Null

# Update the division answers in the table. The table has 20 rows; the 5 rows
# holding data are 1, 5, 9, 13 and 17 (separated by blank spacer rows), each
# with 5 columns. Addressing cells directly (rather than Find/Replace) avoids
# any ambiguity from duplicate "old" values appearing more than once.
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "97÷8=12, 1"
$t.Cell(1, 2).Range.Text = "76÷7=10, 6"
$t.Cell(1, 3).Range.Text = "62÷8=7, 6"
$t.Cell(1, 4).Range.Text = "23÷9=2, 5"
$t.Cell(1, 5).Range.Text = "10÷4=2, 2"

$t.Cell(5, 1).Range.Text = "12÷9=1, 3"
$t.Cell(5, 2).Range.Text = "42÷7=6, 0"
$t.Cell(5, 3).Range.Text = "89÷5=17, 4"
$t.Cell(5, 4).Range.Text = "14÷5=2, 4"
$t.Cell(5, 5).Range.Text = "16÷2=8, 0"

$t.Cell(9, 1).Range.Text = "58÷9=6, 4"
$t.Cell(9, 2).Range.Text = "39÷6=6, 3"
$t.Cell(9, 3).Range.Text = "60÷5=12, 0"
$t.Cell(9, 4).Range.Text = "22÷8=2, 6"
$t.Cell(9, 5).Range.Text = "54÷3=18, 0"

$t.Cell(13, 1).Range.Text = "21÷3=7, 0"
$t.Cell(13, 2).Range.Text = "27÷4=6, 3"
$t.Cell(13, 3).Range.Text = "93÷2=46, 1"
$t.Cell(13, 4).Range.Text = "42÷6=7, 0"
$t.Cell(13, 5).Range.Text = "54÷2=27, 0"

$t.Cell(17, 1).Range.Text = "78÷8=9, 6"
$t.Cell(17, 2).Range.Text = "50÷9=5, 5"
$t.Cell(17, 3).Range.Text = "22÷3=7, 1"
$t.Cell(17, 4).Range.Text = "87÷5=17, 2"
$t.Cell(17, 5).Range.Text = "77÷3=25, 2"
